$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7459603548049927
$ws.Range("B1").Value = 0.5199152231216431
$ws.Range("C1").Value = 0.4039941728115082
$ws.Range("D1").Value = 0.3921604752540588
$ws.Range("E1").Value = 0.4242382347583771
